$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67; this shifts the existing rows 67-143
# down to 68-144 and extends the used range to R144 (matching the
# dimension change from A1:R143 to A1:R144 in the diff).
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new record. The columns
# that are constant for every row in this sheet (A, B, C, E, F, G, I, Q, R)
# are carried over unchanged; the remaining columns (D, H, J, K, L, M, N,
# O, P) hold the new record's values per the diff.
$ws.Range("A67").Value = 4
$ws.Range("B67").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C67").Value = "Los Lagos"
$ws.Range("D67").Value = 44848
$ws.Range("E67").Value = 10
$ws.Range("F67").Value = 100112022
$ws.Range("G67").Value = "Arveja Verde"
$ws.Range("H67").Value = "Perfection"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 40
$ws.Range("K67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("M67").Value = 35000
$ws.Range("N67").Value = "$/malla 25 kilos"
$ws.Range("O67").Value = "Provincia de Huasco"
$ws.Range("P67").Value = 1400
$ws.Range("Q67").Value = 25
$ws.Range("R67").Value = "Hortaliza"
